# also log rollout difficulty proportions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell's formatting (style) from C15 to C16, then set the new date value
$ws.Range("C15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 45427

$ws.Range("A16").Value = "all_training_settings"
$ws.Range("B16").Value = "Virtual"
$ws.Range("E16").Value = "randomEval"
$ws.Range("F16").Value = "random"
$ws.Range("G16").Value = "oncePerTimestep"
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 64
$ws.Range("K16").Value = 256
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 0
$ws.Range("U16").Value = 0
$ws.Range("V16").Value = 0
$ws.Range("W15").Copy($ws.Range("W16"))
$ws.Range("AD16").Value = "check if the training works when all tracks and lights are mixed together"
$ws.Range("AC16").Value = "it learns to complete the easy tracks only???"

# Update sheet view: selection state
$ws.Range("AA21").Select()
